$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current header row (A1:D1): NIS | Nama | NIP | Kode Kelas
# Target header row (A1:E1):  NIS | Nama | NIP Wali | NIP BK | Kode Kelas
# "Kode Kelas" moves from column C to column E, and two new columns
# ("NIP Wali" and "NIP BK") are inserted in its place to support the
# counselor (BK) field added alongside the homeroom teacher (Wali) field.

$kodeKelas = $ws.Range("D1").Value()
$ws.Range("E1").Value = $kodeKelas
$ws.Range("C1").Value = "NIP Wali"
$ws.Range("D1").Value = "NIP BK"

[void]$ws.Range("D2").Select()
